$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 (2025) metrics: total_customers (C6) and new_customers (E6) increase,
# which recalculates new_rate (G6) and returning_rate (H6). retention_rate (F6) unchanged.
$ws.Range("C6").Value = 355
$ws.Range("E6").Value = 73
$ws.Range("G6").Value = 20.56338028169014
$ws.Range("H6").Value = 79.43661971830987
